# Rename the first sheet "INTER_SWITCH_LINKS" -> "SWITCH_TO_SWITCH"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("INTER_SWITCH_LINKS")
$ws1.Name = "SWITCH_TO_SWITCH"

# Make SWITCH_TO_SWITCH the active/selected sheet (moves tabSelected from
# COMPUTE_NODES to SWITCH_TO_SWITCH, and updates the selection on it).
$ws1.Activate()
$ws1.Range("E29").Select()
